# Edit:
#  1) Table on slide 16 switches from the deck's custom table style
#     ({D6F90E39-F23C-4005-939E-EE908EC7093F}) to the built-in style
#     {36F83D91-35A6-4C5F-A382-D75D150776A9}.
#  2) The deck's theme colour scheme (slide master's theme, ppt/theme/
#     theme1.xml) switches from the "Integral" palette to the "Office
#     Theme" palette.

$p = $ppt.ActivePresentation

# --- 1. Table style on slide 16 -----------------------------------------
$s = $p.Slides.Item(16)
$shape = $s.Shapes.Item(3)
$shape.Table.ApplyStyle("{36F83D91-35A6-4C5F-A382-D75D150776A9}")

# --- 2. Switch the deck's theme colours over to the "Office Theme" -----
# palette (was "Integral"). RGB values are packed as R + G*256 + B*65536
# (the VBA long colour order used by ThemeColor.RGB).

$officeColors = @(0, 16777215, 6968388, 15132391, 13998939, 3243501, 10855845, 49407, 12874308, 4697456, 12673797, 7491477)

$masterScheme = $p.SlideMaster.Theme.ThemeColorScheme
for ($i = 1; $i -le $masterScheme.Count; $i++) {
    $masterScheme.Colors($i).RGB = $officeColors[$i - 1]
}
